$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row above row 5 (old row5 -> row6, old row6 -> row7)
$ws.Rows(5).Insert()

# Fill in the map names first (new row 5 needs its other fields too)
$ws.Range("B4").Value = "test맵_1"

# New row 5 (101 / battle_normal variant)
$ws.Range("A5").Value = 101
$ws.Range("B5").Value = "test맵_2"
$ws.Range("C5").Value = "battle_normal"
$ws.Range("F5").Value = "temp_stage_bg"
$ws.Range("G5").Value = "temp_stage_img"
$ws.Range("H5").Value = "map_name_101"
$ws.Range("I5").Value = "map_desc_101"

# Row 6 (was row 5 before insert: 102 / battle_elite): rename map
$ws.Range("B6").Value = "test맵_엘리트"

# Now update the actor_id lists (D column) for all three rows
$ws.Range("D4").Value = "20001,20001"
$ws.Range("D5").Value = "20002,20003"
$ws.Range("D6").Value = "20003,20101,20003"

# Apply text number format to whole actor_id column (D) header + data rows
$ws.Range("D1:D7").NumberFormat = "@"

# Re-add print/page setup (paper size 9 = A4, portrait) lost when resaving
$ws.PageSetup.PaperSize = 9
$ws.PageSetup.Orientation = 1
